# NBAGames.xlsx - "Ran model for 1/24/2020"
# Fills in the "Beat Vegas?" (G) result for the 1/19 games that were
# still pending, and appends the newly-modeled games for 1/24/2021
# (serial date 44220).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Backfill column G ("Beat Vegas?") for rows 111-117 (games on
#    44219 whose outcome is now known).
# ---------------------------------------------------------------
$gResults = @{
    111 = "No"
    112 = "Yes"
    113 = "Yes"
    114 = "Yes"
    115 = "No"
    116 = "Yes"
    117 = "No"
}
foreach ($r in $gResults.Keys) {
    $ws.Cells.Item($r, 7).Value = $gResults[$r]
}

# ---------------------------------------------------------------
# 2) Append the newly modeled games for 1/24/2021 (rows 118-124).
# ---------------------------------------------------------------
$newGames = @(
    @{ Row = 118; Date = 44220; Home = "IND"; Away = "TOR"; Spread = -4.5;  Pred = 1.7;   Diff = -6.2 }
    @{ Row = 119; Date = 44220; Home = "LAC"; Away = "OKC"; Spread = -13;   Pred = -14.5; Diff = 1.5 }
    @{ Row = 120; Date = 44220; Home = "BOS"; Away = "CLE"; Spread = -6;    Pred = -1.9;  Diff = -4.1 }
    @{ Row = 121; Date = 44220; Home = "ORL"; Away = "CHO"; Spread = 1;     Pred = 7.9;   Diff = -6.9 }
    @{ Row = 122; Date = 44220; Home = "MIL"; Away = "ATL"; Spread = -8.5;  Pred = -11.9; Diff = 3.4 }
    @{ Row = 123; Date = 44220; Home = "SAS"; Away = "WAS"; Spread = -7.5;  Pred = 8.9;   Diff = -16.4 }
    @{ Row = 124; Date = 44220; Home = "POR"; Away = "NYK"; Spread = -4;    Pred = -11.1; Diff = 7.1 }
)

foreach ($g in $newGames) {
    $r = $g.Row
    $ws.Cells.Item($r, 1).Value = $g.Date
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat
    $ws.Cells.Item($r, 2).Value = $g.Home
    $ws.Cells.Item($r, 3).Value = $g.Away
    $ws.Cells.Item($r, 4).Value = $g.Spread
    $ws.Cells.Item($r, 5).Value = $g.Pred
    $ws.Cells.Item($r, 6).Value = $g.Diff
}

# ---------------------------------------------------------------
# 3) Cosmetic touch-ups that go along with the refresh: re-fit the
#    date column and restore the working selection/scroll position.
# ---------------------------------------------------------------
$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("A106").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 106
$ws.Range("H11").Select() | Out-Null
